# Generate Report for Handback
#
# Updates the "zh-cn" and "de-de" worksheets of the localization-status
# report with the newly produced handback information for the file
# c683d20c-526e-4ddb-a626-ba759911faa9 (row 8 of each table):
#   - Latest Target File (col I) now links to the handback .md file
#   - Latest Handback File (col J) is filled in with the generated xliff name
#   - Latest Handback DateTime (col K) is updated
#   - Error Detail (col P) reports that the handed-back file is stale
# Also widens the "Error Detail" column so the message is readable.

$wb = $excel.ActiveWorkbook

$handbackMdTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4c231da74728fe8bd5bccb79d1c73b08c591bd43/e2e/c683d20c-526e-4ddb-a626-ba759911faa9.md"
$handbackMdDisplay = "c683d20c-526e-4ddb-a626-ba759911faa9.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/81fc76e0034ff2e0d6e03be1431dbc4949d8b7c5/e2e/c683d20c-526e-4ddb-a626-ba759911faa9.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4c231da74728fe8bd5bccb79d1c73b08c591bd43/e2e/c683d20c-526e-4ddb-a626-ba759911faa9.md."

# Desired column width (OOXML "width" units) for column P; COM ColumnWidth
# is offset from the stored width by the default padding (5/6 chars).
$colWidthTarget = 40 - (5 / 6)

function Update-HandbackRow($ws, $handbackXlfName) {
    # I8: Latest Target File -> hyperlink to the handback markdown file.
    $i8 = $ws.Range("I8")
    $i8.Value = $handbackMdDisplay
    $ws.Hyperlinks.Add($i8, $handbackMdTarget, $null, $null, $handbackMdDisplay) | Out-Null
    $i8.Font.Underline = 2
    $i8.Font.Color = 15570276

    # J8: Latest Handback File -> the generated xliff file name.
    $ws.Range("J8").Value = $handbackXlfName

    # K8: Latest Handback DateTime -> timestamp of the new handback xliff.
    if ($handbackXlfName -eq "c683d20c-526e-4ddb-a626-ba759911faa9.b0b81045a0c6ccaecec6ce3a4805060565956025.zh-cn.xlf") {
        $ws.Range("K8").Value = "2016-09-02 14:53:46"
    } else {
        $ws.Range("K8").Value = "2016-09-02 14:53:54"
    }

    # P8: Error Detail -> report the handback file is not the latest version.
    $ws.Range("P8").Value = $errorDetail

    # Column P (16th column) widened so the message is legible.
    $ws.Range("P1").EntireColumn.ColumnWidth = $colWidthTarget
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Update-HandbackRow $wsZhCn "c683d20c-526e-4ddb-a626-ba759911faa9.b0b81045a0c6ccaecec6ce3a4805060565956025.zh-cn.xlf"

$wsDeDe = $wb.Worksheets.Item("de-de")
Update-HandbackRow $wsDeDe "c683d20c-526e-4ddb-a626-ba759911faa9.b0b81045a0c6ccaecec6ce3a4805060565956025.de-de.xlf"
